# Auto-generated edit script: applies numeric updates per the source diff.
# Each sheet's changed cells are set directly; two cells whose values were
# removed entirely in the diff are cleared (not just zeroed) to match.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1050.2727
$ws.Range("I100").Value = 1050.2727
$ws.Range("K100").Value = 1050.2727
$ws.Range("M100").Value = -509.2727
$ws.Range("H106").Value = 5133974
$ws.Range("I106").Value = 7413518
$ws.Range("K106").Value = 7413518
$ws.Range("M106").Value = -7412887
$ws.Range("H116").Value = 18198.666
$ws.Range("I116").Value = 21655.715
$ws.Range("J116").Value = 6099
$ws.Range("K116").Value = 21655.715
$ws.Range("L116").Value = 6099
$ws.Range("M116").Value = -18213.715
$ws.Range("N116").Value = -12983
$ws.Range("H135").Value = 1342.1666
$ws.Range("I135").Value = 1185.8235
$ws.Range("K135").Value = 10672.4115
$ws.Range("M135").Value = -8137.4115
$ws.Range("H137").Value = 28443.42
$ws.Range("I137").Value = 56909
$ws.Range("J137").Value = 1756.9375
$ws.Range("K137").Value = 170727
$ws.Range("L137").Value = 5270.8125
$ws.Range("M137").Value = -168177
$ws.Range("N137").Value = -10370.8125
$ws.Range("H138").Value = 3132.1829
$ws.Range("I138").Value = 2167.182
$ws.Range("K138").Value = 6501.545999999999
$ws.Range("M138").Value = -1361.545999999999
$ws.Range("H141").Value = 4420.636
$ws.Range("I141").Value = 5515.875
$ws.Range("K141").Value = 16547.625
$ws.Range("M141").Value = -11367.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 76198100
$ws.Range("I32").Value = 50004324
$ws.Range("J32").Value = 111123110
$ws.Range("K32").Value = 50004324
$ws.Range("L32").Value = 111123110
$ws.Range("M32").Value = -50004037
$ws.Range("N32").Value = -111123684
$ws.Range("H33").Value = 6004.8335
$ws.Range("I33").Value = 3200
$ws.Range("K33").Value = 3200
$ws.Range("M33").Value = -2871
$ws.Range("H38").Value = 3650
$ws.Range("I38").Value = 3650
$ws.Range("K38").Value = 3650
$ws.Range("M38").Value = -3183
$ws.Range("H61").Value = 2848.2058
$ws.Range("I61").Value = 2563.4482
$ws.Range("J61").Value = 4499.8
$ws.Range("K61").Value = 2563.4482
$ws.Range("L61").Value = 4499.8
$ws.Range("M61").Value = -2351.4482
$ws.Range("N61").Value = -4923.8
$ws.Range("H63").Value = 178600.62
$ws.Range("I63").Value = 2335
$ws.Range("K63").Value = 2335
$ws.Range("M63").Value = -1649
$ws.Range("H66").Value = 178600.62
$ws.Range("I66").Value = 2335
$ws.Range("K66").Value = 11675
$ws.Range("M66").Value = -8243
$ws.Range("H74").Value = 2030.3636
$ws.Range("I74").Value = 1717.3448
$ws.Range("J74").Value = 4299.75
$ws.Range("K74").Value = 1717.3448
$ws.Range("L74").Value = 4299.75
$ws.Range("M74").Value = -843.3448000000001
$ws.Range("N74").Value = -6047.75
$ws.Range("H77").Value = 2030.3636
$ws.Range("I77").Value = 1717.3448
$ws.Range("J77").Value = 4299.75
$ws.Range("K77").Value = 8586.724
$ws.Range("L77").Value = 21498.75
$ws.Range("M77").Value = -4218.724
$ws.Range("N77").Value = -30234.75
$ws.Range("H110").Value = 1455.2
$ws.Range("I110").Value = 985.75
$ws.Range("J110").Value = 3333
$ws.Range("K110").Value = 985.75
$ws.Range("L110").Value = 3333
$ws.Range("M110").Value = 1059.25
$ws.Range("N110").Value = -7423
$ws.Range("H122").Value = 7284.7856
$ws.Range("I122").Value = 6704.5293
$ws.Range("J122").Value = 8181.5454
$ws.Range("K122").Value = 20113.5879
$ws.Range("L122").Value = 24544.6362
$ws.Range("M122").Value = -17663.5879
$ws.Range("N122").Value = -29444.6362
$ws.Range("H136").Value = 2848.2058
$ws.Range("I136").Value = 2563.4482
$ws.Range("J136").Value = 4499.8
$ws.Range("K136").Value = 7690.344599999999
$ws.Range("L136").Value = 13499.4
$ws.Range("M136").Value = -5140.344599999999
$ws.Range("N136").Value = -18599.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 62008.668
$ws.Range("J51").Value = 62008.668
$ws.Range("L51").Value = 62008.668
$ws.Range("N51").Value = -62990.668
$ws.Range("H102").Value = 27236.3
$ws.Range("I102").Value = 6303
$ws.Range("J102").Value = 110969.5
$ws.Range("K102").Value = 6303
$ws.Range("L102").Value = 110969.5
$ws.Range("M102").Value = -3058
$ws.Range("N102").Value = -117459.5
$ws.Range("H107").Value = 922.05554
$ws.Range("I107").Value = 779.86664
$ws.Range("J107").Value = 1633
$ws.Range("K107").Value = 779.86664
$ws.Range("L107").Value = 1633
$ws.Range("M107").Value = 1140.13336
$ws.Range("N107").Value = -5473

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 194998
$ws.Range("J20").Value = 194998
$ws.Range("L20").Value = 194998
$ws.Range("N20").Value = -195470
$ws.Range("H30").Value = 194998
$ws.Range("J30").Value = 194998
$ws.Range("L30").Value = 194998
$ws.Range("N30").Value = -195180
$ws.Range("H105").Value = 1477.5454
$ws.Range("I105").Value = 917.2222
$ws.Range("K105").Value = 917.2222
$ws.Range("M105").Value = 829.7778
$ws.Range("H128").Value = 194998
$ws.Range("J128").Value = 194998
$ws.Range("L128").Value = 194998
$ws.Range("N128").Value = -204958
$ws.Range("H134").Value = 1909.7812
$ws.Range("I134").Value = 1649.9656
$ws.Range("K134").Value = 4949.8968
$ws.Range("M134").Value = -2414.8968

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 588268.2
$ws.Range("I12").Value = 9.666667
$ws.Range("J12").Value = 909136.4399999999
$ws.Range("K12").Value = 29.000001
$ws.Range("L12").Value = 2727409.32
$ws.Range("M12").Value = 143.999999
$ws.Range("N12").Value = -2727755.32
$ws.Range("H64").Value = 2024.75
$ws.Range("J64").Value = 1550
$ws.Range("L64").Value = 4650
$ws.Range("N64").Value = -5190
$ws.Range("H67").Value = 2024.75
$ws.Range("J67").Value = 1550
$ws.Range("L67").Value = 4650
$ws.Range("N67").Value = -6522
$ws.Range("H103").Value = 159
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 159
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 477
$ws.Range("M103").ClearContents()
$ws.Range("N103").Value = -2235
$ws.Range("H113").Value = 1169.5526
$ws.Range("I113").Value = 1344.6666
$ws.Range("J113").Value = 1115.2069
$ws.Range("K113").Value = 4033.9998
$ws.Range("L113").Value = 3345.620699999999
$ws.Range("M113").Value = -1863.9998
$ws.Range("N113").Value = -7685.620699999999
$ws.Range("H140").Value = 25001824
$ws.Range("I140").Value = 25001824
$ws.Range("K140").Value = 75005472
$ws.Range("M140").Value = -75000292

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3427
$ws.Range("I102").Value = 3355.1
$ws.Range("J102").Value = 3666.6667
$ws.Range("K102").Value = 3355.1
$ws.Range("L102").Value = 3666.6667
$ws.Range("M102").Value = -1733.1
$ws.Range("N102").Value = -6910.6667
$ws.Range("H132").Value = 1547.2413
$ws.Range("I132").Value = 881.8570999999999
$ws.Range("K132").Value = 2645.5713
$ws.Range("M132").Value = -115.5712999999996

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2056.913
$ws.Range("I61").Value = 2368.1667
$ws.Range("J61").Value = 936.4
$ws.Range("K61").Value = 2368.1667
$ws.Range("L61").Value = 936.4
$ws.Range("M61").Value = -2166.1667
$ws.Range("N61").Value = -1340.4
$ws.Range("H113").Value = 2056.913
$ws.Range("I113").Value = 2368.1667
$ws.Range("J113").Value = 936.4
$ws.Range("K113").Value = 2368.1667
$ws.Range("L113").Value = 936.4
$ws.Range("M113").Value = -198.1667000000002
$ws.Range("N113").Value = -5276.4
$ws.Range("H132").Value = 3277.4465
$ws.Range("I132").Value = 2202.2354
$ws.Range("J132").Value = 3746.1282
$ws.Range("K132").Value = 6606.706200000001
$ws.Range("L132").Value = 11238.3846
$ws.Range("M132").Value = -4076.706200000001
$ws.Range("N132").Value = -16298.3846
$ws.Range("H135").Value = 89949.5
$ws.Range("J135").Value = 89949.5
$ws.Range("L135").Value = 89949.5
$ws.Range("N135").Value = -100089.5
$ws.Range("H136").Value = 4704.4165
$ws.Range("J136").Value = 5781.857
$ws.Range("L136").Value = 17345.571
$ws.Range("N136").Value = -22445.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4862.25
$ws.Range("I62").Value = 3350.4
$ws.Range("K62").Value = 3350.4
$ws.Range("M62").Value = -2726.4
$ws.Range("H65").Value = 4862.25
$ws.Range("I65").Value = 3350.4
$ws.Range("K65").Value = 16752
$ws.Range("M65").Value = -13632
$ws.Range("H96").Value = 1990.6666
$ws.Range("I96").Value = 1990.6666
$ws.Range("K96").Value = 1990.6666
$ws.Range("M96").Value = -617.6666
$ws.Range("H108").Value = 99989
$ws.Range("J108").Value = 99989
$ws.Range("L108").Value = 99989
$ws.Range("N108").Value = -107669
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H122").Value = 1913.4546
$ws.Range("I122").Value = 1830.5483
$ws.Range("K122").Value = 5491.644899999999
$ws.Range("M122").Value = -3041.644899999999
$ws.Range("H127").Value = 130000
$ws.Range("J127").Value = 130000
$ws.Range("L127").Value = 130000
$ws.Range("N127").Value = -139920
$ws.Range("H128").Value = 199715
$ws.Range("J128").Value = 199715
$ws.Range("L128").Value = 199715
$ws.Range("N128").Value = -209675
$ws.Range("H132").Value = 2501.4443
$ws.Range("I132").Value = 2019.1951
$ws.Range("K132").Value = 6057.5853
$ws.Range("M132").Value = -3527.5853
$ws.Range("H136").Value = 27498.023
$ws.Range("I136").Value = 1774.079
$ws.Range("J136").Value = 223000
$ws.Range("K136").Value = 5322.237
$ws.Range("L136").Value = 669000
$ws.Range("M136").Value = -2772.237
$ws.Range("N136").Value = -674100
